# fdo#42624: add test for range name import
# Adds range names that reference other range names which are defined
# (and loaded) after the referencing range name, plus a formula cell
# on Sheet2 that exercises the new Global5 -> Global6 -> Sheet2!$B$1 chain.

$wb = $excel.ActiveWorkbook

# New global defined names: Global5 refers to Global6, which in turn is
# only added afterwards - Global5 ends up stored ahead of Global6 in the
# workbook's definedNames list, matching the "name refers to a name that
# is loaded later" import scenario this test targets.
$wb.Names.Add("Global6", "=Sheet2!`$B`$1")
$wb.Names.Add("Global5", "=Global6")

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# New formula cell referencing the freshly added name.
$ws2.Range("A6").Formula = "=Global5"

# Leave the selection state matching a natural "type the formula, hit
# enter" interaction: cursor rests one row below the new entry.
$ws1.Range("A5").Select() | Out-Null
$ws2.Range("A7").Select() | Out-Null
